# Hortaliza, Vega Modelo de Temuco - Apio
# Insert one new weekly price record above the existing row 247, which
# shifts all subsequent records (old rows 247-374) down by one row
# (new rows 248-375). The sheet grows from A1:R374 to A1:R375.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 247, pushing everything below it down one row.
$ws.Rows("247:247").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A247").Value = 10
$ws.Range("B247").Value = "Vega Modelo de Temuco"
$ws.Range("C247").Value = "La Araucanía"
$ws.Range("D247").Value = 44830
$ws.Range("E247").Value = 9
$ws.Range("F247").Value = 100112017
$ws.Range("G247").Value = "Apio"
$ws.Range("H247").Value = "Americana (o)"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 180
$ws.Range("K247").Value = 11000
$ws.Range("L247").Value = 12000
$ws.Range("M247").Value = 11306
$ws.Range("N247").Value = "`$/docena de matas"
$ws.Range("O247").Value = "Provincia del Elquí"
$ws.Range("P247").Value = 1884
$ws.Range("Q247").Value = 6
$ws.Range("R247").Value = "Hortaliza"
